# Updated cryptos list on Sun Apr  2 15:55:17 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the crypto
# listing sheet, row by row, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.234.74'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.804.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.22'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5222'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3819'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07941'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.65'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.356'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.68'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.356'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.806.02'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001091'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06600'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.48'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.970'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.279.19'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.232'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.81'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.50'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.412'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.010.02'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1108'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.064'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.669'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.577'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07233'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.25'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +9.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2180'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.798'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.050'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6217'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.166'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.379'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6054'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.23'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.774'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '126.18'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06826'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.28%  '
